$d = $word.ActiveDocument

# Locate the footnotes touched by this fix dynamically (instead of trusting
# hard-coded collection indices), so the script is robust to the exact
# ordinal position of each footnote:
#   * one footnote ends with a stray "aa" typo that must be dropped
#   * two footnotes are "empty" (their body is just a lone "।") - of those,
#     the very last footnote in the document is a spurious, unreferenced
#     note that should be removed outright, while the other one is simply
#     missing its real text and needs to be filled in.
$footnoteCount = $d.Footnotes.Count

$typoFootnote = $null
$emptyFootnoteIndexes = @()
for ($i = 1; $i -le $footnoteCount; $i++) {
    $txt = $d.Footnotes($i).Range.Text
    if ($txt.EndsWith("aa")) {
        $typoFootnote = $d.Footnotes($i)
    } elseif ($txt.Length -eq 1) {
        $emptyFootnoteIndexes += $i
    }
}

# --- 1. Drop the stray, unreferenced empty footnote (the last one) and the
#        reference + separating space run that point to it in the body.
$lastEmptyIndex = $emptyFootnoteIndexes | Sort-Object | Select-Object -Last 1
$emptyFootnote = $d.Footnotes($lastEmptyIndex)
$emptyFootnote.Delete()

$endPos = $d.Content.End
$trailingSpace = $d.Range($endPos - 2, $endPos - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# --- 2. Fix the footnote with the stray "aa" typo at the end of its text.
$typoFootnote.Range.Text = "བྲིས་པས། སྣར་ཐང་། པེ་ཅིན།"

# --- 3. Fill in the remaining empty footnote's missing text.
$otherEmptyIndex = $emptyFootnoteIndexes | Sort-Object | Select-Object -First 1
$otherEmptyFootnote = $d.Footnotes($otherEmptyIndex)
$otherEmptyFootnote.Range.Text = "རྟོག་པ་རྣམ་རྟོག་པ།_།རྣམ་པར། ཞེས་པར་མ་གཞན་ནང་མེད།"
